$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-07-16 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-07-17 Wednesday", 2)

# Update the multiplication problems in the single table, addressed by
# (row, column) so that duplicate/overlapping values cannot cause mismatches.
$t = $d.Tables.Item(1)

# Row 1
$t.Cell(1, 1).Range.Text = "49×67=3283"
$t.Cell(1, 2).Range.Text = "89×86=7654"
$t.Cell(1, 3).Range.Text = "35×86=3010"
$t.Cell(1, 4).Range.Text = "36×45=1620"
$t.Cell(1, 5).Range.Text = "34×25=850"

# Row 5
$t.Cell(5, 1).Range.Text = "89×79=7031"
$t.Cell(5, 2).Range.Text = "53×30=1590"
$t.Cell(5, 3).Range.Text = "83×73=6059"
$t.Cell(5, 4).Range.Text = "17×21=357"
$t.Cell(5, 5).Range.Text = "21×80=1680"

# Row 10
$t.Cell(10, 1).Range.Text = "76×39=2964"
$t.Cell(10, 2).Range.Text = "65×25=1625"
$t.Cell(10, 3).Range.Text = "20×56=1120"
$t.Cell(10, 4).Range.Text = "96×38=3648"
$t.Cell(10, 5).Range.Text = "88×12=1056"

# Row 15
$t.Cell(15, 1).Range.Text = "97×35=3395"
$t.Cell(15, 2).Range.Text = "31×40=1240"
$t.Cell(15, 3).Range.Text = "64×99=6336"
$t.Cell(15, 4).Range.Text = "60×79=4740"
$t.Cell(15, 5).Range.Text = "95×86=8170"

# Row 20
$t.Cell(20, 1).Range.Text = "91×16=1456"
$t.Cell(20, 2).Range.Text = "54×66=3564"
$t.Cell(20, 3).Range.Text = "83×70=5810"
$t.Cell(20, 4).Range.Text = "47×74=3478"
$t.Cell(20, 5).Range.Text = "91×38=3458"
